$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '39.455.31'
$c.Style = "Normal"
$ws.Range("E2").Value = '  +1.75%  '

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '2.163.05'
$c.Style = "Normal"
$ws.Range("E3").Value = '  +2.82%  '

$ws.Range("E4").Value = '  +0.12%  '

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '228.75'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +0.41%  '

$ws.Range("E6").Value = '  +1.17%  '

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '63.52'
$c.Style = "Normal"
$ws.Range("E7").Value = '  +2.04%  '

$ws.Range("E8").Value = '  +0.05%  '

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.395'
$c.Style = "Normal"
$ws.Range("E9").Value = '  +1.03%  '

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '0.0851'
$c.Style = "Normal"
$ws.Range("E10").Value = '  +1.06%  '

$ws.Range("E11").Value = '  +0.15%  '

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '16.02'
$c.Style = "Normal"
$ws.Range("E12").Value = '  +1.35%  '

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '2.485.66'
$c.Style = "Normal"
$ws.Range("E13").Value = '  +2.89%  '

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '22.00'
$c.Style = "Normal"
$ws.Range("E14").Value = '  -0.45%  '

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '0.812'
$c.Style = "Normal"
$ws.Range("E15").Value = '  +0.19%  '

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '5.52'
$c.Style = "Normal"
$ws.Range("E16").Value = '  -0.13%  '

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '2.170.38'
$c.Style = "Normal"
$ws.Range("E17").Value = '  +2.68%  '

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '39.411.71'
$c.Style = "Normal"
$ws.Range("E18").Value = '  +1.65%  '

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '6.22'
$c.Style = "Normal"
$ws.Range("E19").Value = '  +1.65%  '

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '71.84'
$c.Style = "Normal"
$ws.Range("E20").Value = '  -0.20%  '

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '0.0₃0846'
$c.Style = "Normal"
$ws.Range("E21").Value = '  +0.56%  '

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '229.32'
$c.Style = "Normal"
$ws.Range("E22").Value = '  +0.47%  '

$ws.Range("E23").Value = '  +0.04%  '

$ws.Range("B24").Value = 'Toncoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '2.34'
$c.Style = "Normal"
$ws.Range("E24").Value = '  -0.67%  '

$ws.Range("B25").Value = 'PancakeSwap'
$ws.Range("C25").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '2.36'
$c.Style = "Normal"
$ws.Range("E25").Value = '  +1.67%  '

$ws.Range("B26").Value = 'Cosmos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '9.69'
$c.Style = "Normal"
$ws.Range("E26").Value = '  +0.46%  '

$ws.Range("B27").Value = 'Monero'
$ws.Range("C27").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '172.15'
$c.Style = "Normal"
$ws.Range("E27").Value = '  -0.12%  '

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '0.138'
$c.Style = "Normal"
$ws.Range("E28").Value = '  -1.18%  '

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '19.84'
$c.Style = "Normal"
$ws.Range("E29").Value = '  +2.48%  '

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '1.41'
$c.Style = "Normal"
$ws.Range("E30").Value = '  -0.58%  '

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '2.66'
$c.Style = "Normal"
$ws.Range("E31").Value = '  +4.39%  '

$ws.Range("E32").Value = '  +1.05%  '

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '4.61'
$c.Style = "Normal"
$ws.Range("E33").Value = '  +1.25%  '

$ws.Range("B34").Value = 'InternetComputer(DFINITY)'
$ws.Range("C34").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '4.72'
$c.Style = "Normal"
$ws.Range("E34").Value = '  -1.13%  '

$ws.Range("B35").Value = 'THORChain'
$ws.Range("C35").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '7.03'
$c.Style = "Normal"
$ws.Range("E35").Value = '  -0.22%  '

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '0.0619'
$c.Style = "Normal"
$ws.Range("E36").Value = '  -0.08%  '

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '2.43'
$c.Style = "Normal"
$ws.Range("E37").Value = '  +0.86%  '

$ws.Range("E38").Value = '  +0.00%  '

$ws.Range("E39").Value = '  +0.04%  '

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '102.83'
$c.Style = "Normal"
$ws.Range("E40").Value = '  -0.16%  '

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '0.0228'
$c.Style = "Normal"
$ws.Range("E41").Value = '  -0.14%  '

$ws.Range("E42").Value = '  -1.49%  '

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '1.522.60'
$c.Style = "Normal"
$ws.Range("E43").Value = '  -0.35%  '

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '1.20'
$c.Style = "Normal"
$ws.Range("E44").Value = '  +0.56%  '

$ws.Range("E45").Value = '  +5.19%  '

$ws.Range("E46").Value = '  +0.58%  '

$ws.Range("B47").Value = 'FTXToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '4.28'
$c.Style = "Normal"
$ws.Range("E47").Value = '  +3.53%  '

$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '0.0924'
$c.Style = "Normal"
$ws.Range("E48").Value = '  +1.17%  '

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '7.80'
$c.Style = "Normal"
$ws.Range("E49").Value = '  +0.22%  '

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '2.368.77'
$c.Style = "Normal"
$ws.Range("E50").Value = '  +2.91%  '
